# Updated thesis proposal slides
$p = $ppt.ActivePresentation

# --- 1) Slide 11 ("Poor Support for Performance Anomaly Detection"):
#     split the last bullet into two runs, changing its ending text. ---
$s11 = $p.Slides.Item(11)
$tr11 = $s11.Shapes.Item(2).TextFrame.TextRange
$para11 = $tr11.Paragraphs(6, 1)
$run11 = $para11.Runs(1, 1)
$run11.Text = "Limited support for anomaly detection and "
$run11.InsertAfter("bottleneck identification") | Out-Null

# --- 2) Slide 12 ("Unresolved Issues in the Cloud"):
#     reword the third bullet in place (single run). ---
$s12 = $p.Slides.Item(12)
$tr12 = $s12.Shapes.Item(2).TextFrame.TextRange
$para12 = $tr12.Paragraphs(3, 1)
$run12 = $para12.Runs(1, 1)
$run12.Text = "Difficult to detect performance anomalies and identify bottlenecks"

# --- 3) Insert two new slides before the old "Platform-as-a-Service" run
#     (currently slide 13), both using the same "Title and Content" layout. ---
$layout = $p.Slides.Item(13).CustomLayout

# First insertion -> ends up in position 14 ("Thesis Question").
$sThesis = $p.Slides.AddSlide(13, $layout)
$sThesis.Shapes.Item(1).TextFrame.TextRange.Text = "Thesis Question"
$thesisBody = "Can we enforce design-time governance on web applications developed for a given cloud platform so as to ensure proper versioning, dependency management and conformance to other developer best practices, and also enforce run-time governance on them so as to automatically determine the expected runtime performance of the applications, detect SLA violations and detect performance anomalies and perform root cause analysis, with minimal developer intervention and no invasive instrumentation on the applications?"
$sThesis.Shapes.Item(2).TextFrame.TextRange.Text = $thesisBody

# Second insertion at the same spot pushes the first down -> this one ends
# up in position 13 ("Prelude to Proposal" / cloud platforms governance).
$sPrelude = $p.Slides.AddSlide(13, $layout)
$sPrelude.Shapes.Item(1).TextFrame.TextRange.Text = "Prelude to Proposal"

$preludeLines = @(
    "Automated governance for cloud platforms",
    "Specifying/Learning acceptable operational parameters",
    "Enforcing acceptable operational parameters",
    "Monitoring and detecting deviations from acceptable behavior",
    "Taking corrective/preventive action if necessary"
)
$preludeTr = $sPrelude.Shapes.Item(2).TextFrame.TextRange
$preludeTr.Text = [string]::Join("`r", $preludeLines)
for ($i = 1; $i -le $preludeTr.Paragraphs().Count; $i++) {
    $para = $preludeTr.Paragraphs($i, 1)
    $para.Text = $preludeLines[$i - 1]
    if ($i -ge 2) {
        $para.IndentLevel = 2
    }
}

# --- 4) Delete the old "Prelude to Proposal" / "PaaS clouds" slide, which
#     has now been shifted down to position 18 by the two insertions above. ---
$p.Slides.Item(18).Delete()

# --- 5) Nudge two logo images on the "Aftermath" slide (slide 8). ---
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(4).Left = -127033 / 12700.0
$s8.Shapes.Item(10).Left = 1774338 / 12700.0
